$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: keyword under test changed to a new product name ---
$ws.Range("B3").Value = "Giày NBAL Special Edition"

# --- Row 7: "Tất cả hãng" moved from column D to column C ---
$ws.Range("C7").Value = "Tất cả hãng"
$ws.Range("D7").Value = ""

# --- Result / detail columns (G = Detail message, H = Pass/Fail) for rows 2-9 ---
# entered in row order, column G before H, matching how the shared-string
# table grows in the authored workbook
$ws.Range("G2").Value = "Thành công: 8 SP"
$ws.Range("H2").Value = "PASS"

$ws.Range("G3").Value = "Lỗi: Không tìm thấy sản phẩm"
$ws.Range("H3").Value = "FAIL"

$ws.Range("G4").Value = "Đúng: 0 kết quả"
$ws.Range("H4").Value = "PASS"

$ws.Range("G5").Value = "Thành công: 8 SP"
$ws.Range("H5").Value = "PASS"

$ws.Range("G6").Value = "Thành công: 8 SP"
$ws.Range("H6").Value = "PASS"

$ws.Range("G7").Value = "Lỗi hệ thống: Cannot locate option with text: Tất cả hãng`nFor documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception`nBuild info: version: '4.14.1', revision: '03f8ede370'`nSystem info: os.name: 'Windows 11', os.arch: 'amd64', os.version: '10.0', java.version: '17.0.10'`nDriver info: driver.version: unknown"
$ws.Range("H7").Value = "FAIL"
# long wrapped text in G7 otherwise triggers auto row-height growth; pin it back
$ws.Rows(7).RowHeight = 15.75

$ws.Range("G8").Value = "Đúng: 0 kết quả"
$ws.Range("H8").Value = "PASS"

$ws.Range("G9").Value = "Thành công: 8 SP"
$ws.Range("H9").Value = "PASS"

# --- Column widths: re-affirm the explicit widths (pixel-exact) ---
$ws.Columns("D").Width = 140
$ws.Columns("F").Width = 220
$ws.Columns("G").Width = 290

# --- Selection moved to B3 ---
$ws.Range("B3").Select()
